$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reduce some light requirements (min_skylight, column N) to 0 for several plants,
# and bump the fertilize value (column B) for two rows.

$ws.Range("N9").Value = 0

$ws.Range("B19").Value = 20
$ws.Range("N19").Value = 0

$ws.Range("B20").Value = 20
$ws.Range("N20").Value = 0

$ws.Range("N62").Value = 0

$ws.Range("N75").Value = 0
$ws.Range("N76").Value = 0
$ws.Range("N78").Value = 0
$ws.Range("N79").Value = 0
$ws.Range("N81").Value = 0
$ws.Range("N82").Value = 0

$ws.Range("N86").Value = 0
$ws.Range("N87").Value = 0
$ws.Range("N88").Value = 0
$ws.Range("N89").Value = 0
$ws.Range("N90").Value = 0
$ws.Range("N91").Value = 0

# Leave the cursor/selection where the author's last click landed.
$ws.Range("N9").Select()
